$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: two new "Pera" price observations (Packham's Triumph, Primera
# and Segunda) were reported for Feria Lagunitas de Puerto Montt. They belong
# right after the existing row for that market (row 290 in the old layout), so
# insert two blank rows there; everything below shifts down by two rows.
$ws.Rows("290:291").Insert()

function Set-RowValues($rowIndex, $values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
    }
}

# New row 290: Packham's Triumph, Primera
Set-RowValues 290 @(
    4,
    "Feria Lagunitas de Puerto Montt",
    "Los Lagos",
    44783,
    10,
    "Fruta",
    100104,
    "Frutos de pepita",
    100104005,
    "Pera",
    "Packham's Triumph",
    "Primera",
    200,
    15000,
    16000,
    15500,
    "$/caja 15 kilos empedrada",
    "Región de O'Higgins",
    1033,
    15
)

# New row 291: Packham's Triumph, Segunda
Set-RowValues 291 @(
    4,
    "Feria Lagunitas de Puerto Montt",
    "Los Lagos",
    44783,
    10,
    "Fruta",
    100104,
    "Frutos de pepita",
    100104005,
    "Pera",
    "Packham's Triumph",
    "Segunda",
    100,
    13000,
    13000,
    13000,
    "$/caja 15 kilos empedrada",
    "Región de O'Higgins",
    867,
    15
)
